# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" column (D) for every localized-file
# row that is about to be (re-)handed off: status "Handback transform failed"
# (row 4) and status "Ready for handoff" (rows 6-10), on both the "zh-cn"
# and "de-de" language sheets. Rows already "Handed back: in sync with en-US"
# (rows 2-3), "In Translation" (row 5), or the non-localized config row (11)
# are left untouched.

$wb = $excel.ActiveWorkbook

$rowsToStamp = @(4, 6, 7, 8, 9, 10)

$zhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToStamp) {
    $zhCn.Cells.Item($r, 4).Value = "2016-02-18 10:42:08"
}

$deDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToStamp) {
    $deDe.Cells.Item($r, 4).Value = "2016-02-18 10:42:21"
}
